$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7602.769
$ws.Range("I32").Value = 6424.8335
$ws.Range("J32").Value = 8612.429
$ws.Range("K32").Value = 6424.8335
$ws.Range("L32").Value = 8612.429
$ws.Range("M32").Value = -6098.8335
$ws.Range("N32").Value = -9264.429
$ws.Range("H127").Value = 1277
$ws.Range("I127").Value = 1277
$ws.Range("K127").Value = 3831
$ws.Range("M127").Value = 1129
$ws.Range("H137").Value = 4093.5334
$ws.Range("I137").Value = 2479.111
$ws.Range("J137").Value = 6515.1665
$ws.Range("K137").Value = 7437.333
$ws.Range("L137").Value = 19545.4995
$ws.Range("M137").Value = -4887.333
$ws.Range("N137").Value = -24645.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3392.4443
$ws.Range("I2").Value = 2361
$ws.Range("K2").Value = 2361
$ws.Range("M2").Value = -2248
$ws.Range("H45").Value = 1786.5714
$ws.Range("I45").Value = 1666.2
$ws.Range("J45").Value = 2087.5
$ws.Range("K45").Value = 1666.2
$ws.Range("L45").Value = 2087.5
$ws.Range("M45").Value = -1289.2
$ws.Range("N45").Value = -2841.5
$ws.Range("H61").Value = 2455.4849
$ws.Range("I61").Value = 2452.6775
$ws.Range("K61").Value = 2452.6775
$ws.Range("M61").Value = -2240.6775
$ws.Range("H116").Value = 3392.4443
$ws.Range("I116").Value = 2361
$ws.Range("K116").Value = 2361
$ws.Range("M116").Value = -67
$ws.Range("H132").Value = 19236090
$ws.Range("I132").Value = 3902.5908
$ws.Range("J132").Value = 125013120
$ws.Range("K132").Value = 11707.7724
$ws.Range("L132").Value = 375039360
$ws.Range("M132").Value = -9177.7724
$ws.Range("N132").Value = -375044420
$ws.Range("H136").Value = 2455.4849
$ws.Range("I136").Value = 2452.6775
$ws.Range("K136").Value = 7358.032499999999
$ws.Range("M136").Value = -4808.032499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3392.4443
$ws.Range("I3").Value = 2361
$ws.Range("K3").Value = 2361
$ws.Range("M3").Value = -2247
$ws.Range("H20").Value = 2633.7222
$ws.Range("I20").Value = 1568.75
$ws.Range("J20").Value = 4763.6665
$ws.Range("K20").Value = 1568.75
$ws.Range("L20").Value = 4763.6665
$ws.Range("M20").Value = -1321.75
$ws.Range("N20").Value = -5257.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1163.3334
$ws.Range("I16").Value = 1163.3334
$ws.Range("K16").Value = 1163.3334
$ws.Range("M16").Value = -876.3334
$ws.Range("H31").Value = 1845.6666
$ws.Range("I31").Value = 1159
$ws.Range("K31").Value = 1159
$ws.Range("M31").Value = -864
$ws.Range("H34").Value = 1845.6666
$ws.Range("I34").Value = 1159
$ws.Range("K34").Value = 1159
$ws.Range("M34").Value = -957
$ws.Range("H62").Value = 17957.125
$ws.Range("I62").Value = 12615.556
$ws.Range("K62").Value = 12615.556
$ws.Range("M62").Value = -11991.556
$ws.Range("H65").Value = 17957.125
$ws.Range("I65").Value = 12615.556
$ws.Range("K65").Value = 63077.78
$ws.Range("M65").Value = -59957.78
$ws.Range("H113").Value = 1163.3334
$ws.Range("I113").Value = 1163.3334
$ws.Range("K113").Value = 1163.3334
$ws.Range("M113").Value = 1006.6666
$ws.Range("H135").Value = 77626.336
$ws.Range("J135").Value = 77626.336
$ws.Range("L135").Value = 77626.336
$ws.Range("N135").Value = -87766.336
$ws.Range("H140").Value = 18000
$ws.Range("J140").Value = 18000
$ws.Range("L140").Value = 18000
$ws.Range("N140").Value = -28360
$ws.Range("H141").Value = 70662
$ws.Range("J141").Value = 70662
$ws.Range("L141").Value = 70662
$ws.Range("N141").Value = -81022

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3249.6875
$ws.Range("J69").Value = 3249.6875
$ws.Range("L69").Value = 9749.0625
$ws.Range("N69").Value = -11371.0625
$ws.Range("H72").Value = 3249.6875
$ws.Range("J72").Value = 3249.6875
$ws.Range("L72").Value = 29247.1875
$ws.Range("N72").Value = -37359.1875
$ws.Range("H88").Value = 2432.3333
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 2432.3333
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H116").Value = 46234.59
$ws.Range("I116").Value = 98641.14
$ws.Range("J116").Value = 9550
$ws.Range("K116").Value = 295923.42
$ws.Range("L116").Value = 28650
$ws.Range("M116").Value = -292481.42
$ws.Range("N116").Value = -35534
$ws.Range("H137").Value = 539911.2
$ws.Range("J137").Value = 837428.7
$ws.Range("L137").Value = 2512286.1
$ws.Range("N137").Value = -2522486.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 7000
$ws.Range("I99").Value = 7000
$ws.Range("K99").Value = 7000
$ws.Range("M99").Value = -4754
$ws.Range("H122").Value = 2374.1924
$ws.Range("I122").Value = 2531.2
$ws.Range("K122").Value = 7593.599999999999
$ws.Range("M122").Value = -5143.599999999999
$ws.Range("H132").Value = 2798.6667
$ws.Range("I132").Value = 2748.5833
$ws.Range("K132").Value = 8245.749899999999
$ws.Range("M132").Value = -5715.749899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3036.1667
$ws.Range("I7").Value = 3243.4
$ws.Range("K7").Value = 3243.4
$ws.Range("M7").Value = -3131.4
$ws.Range("H109").Value = 25189.334
$ws.Range("J109").Value = 25189.334
$ws.Range("L109").Value = 25189.334
$ws.Range("N109").Value = -27963.334
$ws.Range("H126").Value = 3036.1667
$ws.Range("I126").Value = 3243.4
$ws.Range("K126").Value = 9730.200000000001
$ws.Range("M126").Value = -7260.200000000001
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("H136").Value = 4733.793
$ws.Range("I136").Value = 3388.7368
$ws.Range("K136").Value = 10166.2104
$ws.Range("M136").Value = -7616.2104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 30083.5
$ws.Range("J53").Value = 30083.5
$ws.Range("L53").Value = 30083.5
$ws.Range("N53").Value = -31297.5
$ws.Range("H70").Value = 27958.947
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 27958.947
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 27958.947
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -28588.947
$ws.Range("H73").Value = 27958.947
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 27958.947
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 27958.947
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -30142.947
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H122").Value = 4207.533
$ws.Range("I122").Value = 5136.3335
$ws.Range("K122").Value = 15409.0005
$ws.Range("M122").Value = -12959.0005
$ws.Range("H126").Value = 3175.647
$ws.Range("I126").Value = 2539
$ws.Range("J126").Value = 3440.9167
$ws.Range("K126").Value = 7617
$ws.Range("L126").Value = 10322.7501
$ws.Range("M126").Value = -5147
$ws.Range("N126").Value = -15262.7501
$ws.Range("H132").Value = 2676.6
$ws.Range("I132").Value = 2676.6
$ws.Range("K132").Value = 8029.799999999999
$ws.Range("M132").Value = -5499.799999999999
$ws.Range("H133").Value = 50711
$ws.Range("J133").Value = 50711
$ws.Range("L133").Value = 50711
$ws.Range("N133").Value = -60831
$ws.Range("H136").Value = 831.4
$ws.Range("I136").Value = 831.4
$ws.Range("K136").Value = 2494.2
$ws.Range("M136").Value = 55.80000000000018
